$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 95

$ws.Cells.Item($row, 1).Value = "Kindergarden"
$ws.Cells.Item($row, 2).Value = "Kindergarden Den Haag Appelgaard"
$ws.Cells.Item($row, 3).Value = "KDV"

# The report date looks like a date literal; force it to stay as plain text
# (matching the rest of the sheet's "Rapportdatum" column) and drop any
# number-format styling that gets attached when typing a date-looking value.
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "2024-07-29"
$ws.Cells.Item($row, 4).ClearFormats()

$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 1
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
